# Word Update - Version 2 (V2 word changes)
#
# 1. Merge the two runs of paragraph 3 ("We want to co" / "ntrol these too",
#    split by the hidden "_GoBack" bookmark) into a single run reading
#    "We want to control these too", removing that bookmark.
# 2. Insert, after that paragraph: a blank paragraph, a bold+underlined
#    "Version 2" paragraph, and a blank bold+underlined paragraph.
# 3. Append a final paragraph "This is some more text", re-creating the
#    "_GoBack" bookmark between "This is som" and "e more text" so the
#    run split matches the original straddling pattern.

$d = $word.ActiveDocument

# --- Step 1: drop the stale hidden bookmark, then merge/normalize the text ---
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

$null = $d.Content.Find.Execute(
    "We want to co" + "ntrol these too",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "We want to control these too", 2)

# --- Step 2: build the four new paragraphs after "We want to control..." ---
$pControl = $d.Paragraphs.Item(3)
$pControl.Range.InsertParagraphAfter()      # new blank paragraph

$pBlank1 = $d.Paragraphs.Item(4)
$pBlank1.Range.InsertParagraphAfter()       # paragraph that becomes "Version 2"

$pVersion = $d.Paragraphs.Item(5)
$pVersion.Range.InsertParagraphAfter()      # blank bold/underlined paragraph

$pBlank2 = $d.Paragraphs.Item(6)
$pBlank2.Range.InsertParagraphAfter()       # paragraph that becomes the final text

$pFinal = $d.Paragraphs.Item(7)

$pVersion.Range.Text = "Version 2"
$pVersion.Range.Font.Bold = 1
$pVersion.Range.Font.Underline = 1

$pBlank2.Range.Font.Bold = 1
$pBlank2.Range.Font.Underline = 1

$pFinal.Range.Text = "This is some more text"

# --- Step 3: split the final paragraph's text with a re-created "_GoBack" ---
$splitPos = $pFinal.Range.Start + "This is som".Length
$bmRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
